$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.362.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.604.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("E10").Value = '  +1.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0857'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.831.27'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.609.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.371.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.00%  '
$ws.Range("E19").Value = '  +4.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0723'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.40'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.83%  '
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("B33").Value = 'Maker'
$ws.Range("C33").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.484.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.28%  '
$ws.Range("E34").Value = '  -1.25%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.559'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0164'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.819'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E42").Value = '  +1.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.937'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.742.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.82%  '
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.83%  '
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0959'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("E51").Value = '  +0.10%  '
